$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.367.17"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.686.27"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +3.83%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "408.95"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.19"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.76%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.678.55"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.625"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.40%  "
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.731"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -5.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.165"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -9.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000335"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -6.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.20"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.96"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.295.31"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.35%  "
$ws.Range("E16").Value = "  -1.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.685.92"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.03"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.99"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.10%  "
$ws.Range("E20").Value = "  -3.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "65.464.15"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "421.31"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -5.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.10"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +15.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.52"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.02"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -5.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "36.32"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +5.76%  "
$ws.Range("E27").Value = "  -5.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.50"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.11"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.54"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.119"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.89%  "
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.71"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.54%  "
$ws.Range("E33").Value = "  -4.41%  "
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "40.86"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +5.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.81"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0468"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.58%  "
$ws.Range("E39").Value = "  +24.67%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.140"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.54%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.996"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₃0654"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -19.36%  "
$ws.Range("B43").Value = "LidoDAOToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.35"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.02"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +25.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "143.00"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.70%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.06"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.27%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.26"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.53%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.04"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +17.27%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.79"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -7.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.52"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -8.98%  "
$ws.Range("E51").Value = "  -5.47%  "
